# Updated cryptos list on Thu Oct 12 09:43:32 UTC 2023 with GitHub Actions
#
# Refresh the price/1h-volume-change snapshot in the cryptos table with the
# latest values from coinranking.com. Two rows (Solana/Cardano) also swap
# rank position, so their Coin/Link/Price/Volume cells all move together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.769.87'
$ws.Range('E2').Value = '  -1.52%  '
$ws.Range('D3').Value = '1.552.21'
$ws.Range('E3').Value = '  -1.34%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '204.78'
$ws.Range('E5').Value = '  -1.12%  '
$ws.Range('E6').Value = '  -1.91%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = '0.246'
$ws.Range('E8').Value = '  -0.75%  '
$ws.Range('B9').Value = 'Solana'
$ws.Range('C9').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D9').Value = '21.43'
$ws.Range('E9').Value = '  -4.01%  '
$ws.Range('D10').Value = '0.0582'
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('D12').Value = '1.778.93'
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('D13').Value = '1.563.10'
$ws.Range('E13').Value = '  -0.53%  '
$ws.Range('D14').Value = '3.69'
$ws.Range('E14').Value = '  -2.14%  '
# Preserve the trailing zero as literal text (would otherwise be
# normalised away by numeric auto-detection on assignment).
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.510'
$ws.Range('E15').Value = '  -1.85%  '
$ws.Range('D16').Value = '26.780.78'
$ws.Range('E16').Value = '  -1.46%  '
$ws.Range('D17').Value = '61.15'
$ws.Range('E17').Value = '  -2.41%  '
$ws.Range('D18').Value = '214.13'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').Value = '7.27'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').Value = '0.0₃0681'
$ws.Range('E20').Value = '  -0.82%  '
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('D22').Value = '4.08'
$ws.Range('E22').Value = '  -1.42%  '
$ws.Range('D23').Value = '9.08'
$ws.Range('E23').Value = '  -3.46%  '
$ws.Range('D24').Value = '2.01'
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').Value = '152.02'
$ws.Range('E25').Value = '  -0.46%  '
$ws.Range('D26').Value = '6.53'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('D27').Value = '14.89'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  -2.23%  '
$ws.Range('E30').Value = '  -0.34%  '
$ws.Range('E31').Value = '  -3.05%  '
$ws.Range('E32').Value = '  -0.42%  '
$ws.Range('D33').Value = '1.366.28'
$ws.Range('E33').Value = '  -2.34%  '
# Preserve the trailing zero as literal text (would otherwise be
# normalised away by numeric auto-detection on assignment).
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.90'
$ws.Range('E34').Value = '  -0.64%  '
$ws.Range('E35').Value = '  -4.00%  '
$ws.Range('E36').Value = '  -0.43%  '
$ws.Range('D37').Value = '0.921'
$ws.Range('E37').Value = '  -2.67%  '
$ws.Range('D38').Value = '0.0162'
$ws.Range('E38').Value = '  -2.14%  '
$ws.Range('D39').Value = '0.521'
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('D40').Value = '0.801'
$ws.Range('E40').Value = '  -1.91%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  +3.99%  '
$ws.Range('E43').Value = '  -0.35%  '
$ws.Range('E44').Value = '  +0.22%  '
$ws.Range('D45').Value = '1.76'
$ws.Range('E45').Value = '  -2.79%  '
$ws.Range('D46').Value = '63.03'
$ws.Range('E46').Value = '  -1.16%  '
$ws.Range('E47').Value = '  -2.61%  '
$ws.Range('D48').Value = '1.690.56'
$ws.Range('E48').Value = '  -0.94%  '
$ws.Range('D49').Value = '86.14'
$ws.Range('E49').Value = '  +0.17%  '
$ws.Range('D50').Value = '0.0516'
$ws.Range('E50').Value = '  +4.54%  '
$ws.Range('D51').Value = '0.0₇0980'
$ws.Range('E51').Value = '  -0.14%  '
